# Renamed few transcripts. Updated the DataSheet.
# Column D ("Speaker"): "Davis" -> "T", "Student" -> "S" for the specified rows.
# Column F ("Teacher Tag"): "3 - getting students to relate" -> "3 - getting SS to relate" for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where Speaker was "Davis" and should become "T"
$davisRows = @(12,16,27,28,29,30,37,47,51,53,62,63,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,103,109,110,114,128,141,150,151,168,169,172,173)

# Rows where Speaker was "Student" and should become "S"
$studentRows = @(18,19,32,36,40,42,44,46,52,59,105,107,115,117,119,127,129,132,133,134,135,136,137,142,145,146,152,154,155,157,160,175)

foreach ($r in $davisRows) {
    $ws.Range("D$r").Value = "T"
}

foreach ($r in $studentRows) {
    $ws.Range("D$r").Value = "S"
}

# Rows where Teacher Tag text mentions "students" and should become "SS"
$tagRows = @(62,85,87,173)

foreach ($r in $tagRows) {
    $ws.Range("F$r").Value = "3 - getting SS to relate"
}
